$d = $word.ActiveDocument

# 1. Remove "Python, " from the programming-languages list.
$d.Content.Find.Execute(
    "cript, TypeScript, HTML/CSS, Python, C++, C#, ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "cript, TypeScript, HTML/CSS, C++, C#, ", 2
) | Out-Null

# 2. Drop the ", PHP" entry entirely (list becomes ...SQL, Rust).
$d.Content.Find.Execute(
    "SQL, PHP, Rust",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "SQL, Rust", 2
) | Out-Null

# 3. Add a new bullet right after "...attended standup meetings", reusing
#    the same list numbering (numId 15) as its sibling bullets.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*standup meetings*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ge 1) {
    # Create a fresh blank paragraph right after the target bullet, then
    # fill it via raw OOXML so the new list item's numPr/rPr match the
    # existing bullets exactly (numId 15, Times New Roman 12pt).
    $target = $d.Paragraphs.Item($targetIndex)
    $target.Range.InsertParagraphAfter() | Out-Null
    $newPara = $d.Paragraphs.Item($targetIndex + 1)
    $insertionRange = $newPara.Range
    $insertionRange.Collapse(1)
    $xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="15"/></w:numPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>Debugged, fixed and maintained 60% of the codebase.</w:t></w:r></w:p>
'@
    $insertionRange.InsertXML($xml)
}
